$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.98769999999999
$ws.Range("B21").Value = 5.667399999999994
$ws.Range("B23").Value = 5.723899999999999
$ws.Range("B25").Value = 5.884699999999994
